# B1--and-B2-PowerPoint.pptx edit
#
# 1) The table on slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS") gets a new
#    built-in table style applied from the Table Design gallery
#    ({7F0F8574-2EF2-4012-8597-B37BF814D354} -> {3103B075-8CB8-413A-9B73-E85E3D678D55}).
#
# 2) The deck's theme ("Integral" / Red-Violet) is switched back to the
#    default Office theme - i.e. the slide master's 12 theme colors are
#    reset to the standard Office palette.

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{3103B075-8CB8-413A-9B73-E85E3D678D55}")

# --- 2) Theme colors: Red Violet / Integral -> default Office -------------
# ThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeHex = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeHex.Count; $i++) {
    $hex = $officeHex[$i - 1]
    $rr = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $gg = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $bb = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $rgb = $rr + ($gg * 256) + ($bb * 65536)
    $themeColors.Colors($i).RGB = $rgb
}
